# Auto-generated edit script: refresh market-price/profit columns (H-N)
# across all 8 crafting-job sheets, per the scheduled Hyperion price-update run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 4789.716
$ws.Range("J17").Value = 4789.716
$ws.Range("L17").Value = 14369.148
$ws.Range("N17").Value = -14705.148
$ws.Range("H33").Value = 617.2727
$ws.Range("I33").Value = 682.41174
$ws.Range("J33").Value = 395.8
$ws.Range("K33").Value = 682.41174
$ws.Range("L33").Value = 395.8
$ws.Range("M33").Value = -453.41174
$ws.Range("N33").Value = -853.8
$ws.Range("H43").Value = 1685.7778
$ws.Range("I43").Value = 1767.5
$ws.Range("K43").Value = 1767.5
$ws.Range("M43").Value = -1698.5
$ws.Range("H62").Value = 6480.1333
$ws.Range("I62").Value = 3477.077
$ws.Range("K62").Value = 3477.077
$ws.Range("M62").Value = -2853.077
$ws.Range("H65").Value = 6480.1333
$ws.Range("I65").Value = 3477.077
$ws.Range("K65").Value = 17385.385
$ws.Range("M65").Value = -14265.385
$ws.Range("H74").Value = 7386.731
$ws.Range("I74").Value = 4098.5
$ws.Range("J74").Value = 7660.75
$ws.Range("K74").Value = 4098.5
$ws.Range("L74").Value = 7660.75
$ws.Range("M74").Value = -3162.5
$ws.Range("N74").Value = -9532.75
$ws.Range("H77").Value = 7386.731
$ws.Range("I77").Value = 4098.5
$ws.Range("J77").Value = 7660.75
$ws.Range("K77").Value = 20492.5
$ws.Range("L77").Value = 38303.75
$ws.Range("M77").Value = -15812.5
$ws.Range("N77").Value = -47663.75
$ws.Range("H88").Value = 3109.1333
$ws.Range("J88").Value = 3338.4783
$ws.Range("L88").Value = 3338.4783
$ws.Range("N88").Value = -4150.478300000001
$ws.Range("H91").Value = 3109.1333
$ws.Range("J91").Value = 3338.4783
$ws.Range("L91").Value = 3338.4783
$ws.Range("N91").Value = -6146.478300000001
$ws.Range("H92").Value = 620.08105
$ws.Range("I92").Value = 696.40625
$ws.Range("K92").Value = 696.40625
$ws.Range("M92").Value = 551.59375
$ws.Range("H97").Value = 2703.3333
$ws.Range("J97").Value = 3055
$ws.Range("L97").Value = 9165
$ws.Range("N97").Value = -10157
$ws.Range("H104").Value = 294.1875
$ws.Range("I104").Value = 305.8
$ws.Range("K104").Value = 917.4000000000001
$ws.Range("M104").Value = 829.5999999999999
$ws.Range("H115").Value = 483.5
$ws.Range("I115").Value = 483.5
$ws.Range("K115").Value = 1450.5
$ws.Range("M115").Value = 116.5
$ws.Range("H125").Value = 3458.738
$ws.Range("J125").Value = 3631.9211
$ws.Range("L125").Value = 32687.2899
$ws.Range("N125").Value = -37607.2899
$ws.Range("H129").Value = 29412652
$ws.Range("I129").Value = 31250756
$ws.Range("K129").Value = 93752268
$ws.Range("M129").Value = -93747268
$ws.Range("H132").Value = 1415.8914
$ws.Range("I132").Value = 1415.8914
$ws.Range("K132").Value = 4247.674199999999
$ws.Range("M132").Value = -1717.674199999999
$ws.Range("H137").Value = 3002.8206
$ws.Range("I137").Value = 2758
$ws.Range("J137").Value = 3212.6667
$ws.Range("K137").Value = 8274
$ws.Range("L137").Value = 9638.000100000001
$ws.Range("M137").Value = -5724
$ws.Range("N137").Value = -14738.0001
$ws.Range("H141").Value = 4145.5713
$ws.Range("I141").Value = 4355.1035
$ws.Range("J141").Value = 3132.8333
$ws.Range("K141").Value = 13065.3105
$ws.Range("L141").Value = 9398.499899999999
$ws.Range("M141").Value = -7885.3105
$ws.Range("N141").Value = -19758.4999

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4173.196
$ws.Range("I32").Value = 2748.6575
$ws.Range("K32").Value = 2748.6575
$ws.Range("M32").Value = -2461.6575
$ws.Range("H61").Value = 1571
$ws.Range("I61").Value = 1571
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1571
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1359
$ws.Range("H74").Value = 75485.14
$ws.Range("I74").Value = 60472.4
$ws.Range("K74").Value = 60472.4
$ws.Range("M74").Value = -59598.4
$ws.Range("H77").Value = 75485.14
$ws.Range("I77").Value = 60472.4
$ws.Range("K77").Value = 302362
$ws.Range("M77").Value = -297994
$ws.Range("H97").Value = 18560.25
$ws.Range("I97").Value = 16967.545
$ws.Range("K97").Value = 16967.545
$ws.Range("M97").Value = -16471.545
$ws.Range("H102").Value = 5958.7144
$ws.Range("I102").Value = 5778
$ws.Range("J102").Value = 6726.75
$ws.Range("K102").Value = 5778
$ws.Range("L102").Value = 6726.75
$ws.Range("M102").Value = -4156
$ws.Range("N102").Value = -9970.75
$ws.Range("H110").Value = 1134
$ws.Range("I110").Value = 1150.3529
$ws.Range("K110").Value = 1150.3529
$ws.Range("M110").Value = 894.6470999999999
$ws.Range("H124").Value = 7725.8
$ws.Range("J124").Value = 7725.8
$ws.Range("L124").Value = 7725.8
$ws.Range("N124").Value = -17545.8
$ws.Range("H132").Value = 2413.2964
$ws.Range("I132").Value = 1845.9048
$ws.Range("K132").Value = 5537.7144
$ws.Range("M132").Value = -3007.7144
$ws.Range("H136").Value = 1571
$ws.Range("I136").Value = 1571
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4713
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -2163

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value = 1795.5
$ws.Range("J20").Value = 2749
$ws.Range("L20").Value = 2749
$ws.Range("N20").Value = -3243
$ws.Range("H86").Value = 12757.368
$ws.Range("I86").Value = 10381.857
$ws.Range("J86").Value = 19408.8
$ws.Range("K86").Value = 10381.857
$ws.Range("L86").Value = 19408.8
$ws.Range("M86").Value = -9258.857
$ws.Range("N86").Value = -21654.8
$ws.Range("H89").Value = 12757.368
$ws.Range("I89").Value = 10381.857
$ws.Range("J89").Value = 19408.8
$ws.Range("K89").Value = 51909.285
$ws.Range("L89").Value = 97044
$ws.Range("M89").Value = -46293.285
$ws.Range("N89").Value = -108276
$ws.Range("H94").Value = 5117.6
$ws.Range("I94").Value = 761.0714
$ws.Range("K94").Value = 761.0714
$ws.Range("M94").Value = -310.0714
$ws.Range("H99").Value = 48691.184
$ws.Range("I99").Value = 73730.93
$ws.Range("K99").Value = 73730.93
$ws.Range("M99").Value = -72232.93
$ws.Range("H132").Value = 90000
$ws.Range("J132").Value = 90000
$ws.Range("L132").Value = 90000
$ws.Range("N132").Value = -100120
$ws.Range("H134").Value = 2782.25
$ws.Range("I134").Value = 848.9667
$ws.Range("J134").Value = 8582.1
$ws.Range("K134").Value = 2546.9001
$ws.Range("L134").Value = 25746.3
$ws.Range("M134").Value = -11.90009999999984
$ws.Range("N134").Value = -30816.3

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H6").Value = 621.75
$ws.Range("I6").Value = 787.3333
$ws.Range("J6").Value = 125
$ws.Range("K6").Value = 787.3333
$ws.Range("L6").Value = 125
$ws.Range("M6").Value = -674.3333
$ws.Range("N6").Value = -351
$ws.Range("H31").Value = 4081.742
$ws.Range("I31").Value = 2716.5
$ws.Range("J31").Value = 5538
$ws.Range("K31").Value = 2716.5
$ws.Range("L31").Value = 5538
$ws.Range("M31").Value = -2421.5
$ws.Range("N31").Value = -6128
$ws.Range("H34").Value = 4081.742
$ws.Range("I34").Value = 2716.5
$ws.Range("J34").Value = 5538
$ws.Range("K34").Value = 2716.5
$ws.Range("L34").Value = 5538
$ws.Range("M34").Value = -2514.5
$ws.Range("N34").Value = -5942
$ws.Range("H58").Value = 1697.7333
$ws.Range("I58").Value = 1462.2142
$ws.Range("J58").Value = 4995
$ws.Range("K58").Value = 1462.2142
$ws.Range("L58").Value = 4995
$ws.Range("M58").Value = -1259.2142
$ws.Range("N58").Value = -5401
$ws.Range("H105").Value = 4724.75
$ws.Range("I105").Value = 4000
$ws.Range("J105").Value = 4966.3335
$ws.Range("K105").Value = 4000
$ws.Range("L105").Value = 4966.3335
$ws.Range("M105").Value = -2253
$ws.Range("N105").Value = -8460.3335
$ws.Range("H107").Value = 3837.75
$ws.Range("I107").Value = 3948.3333
$ws.Range("K107").Value = 3948.3333
$ws.Range("M107").Value = -2028.3333
$ws.Range("H122").Value = 3523.9546
$ws.Range("I122").Value = 2765.5386
$ws.Range("K122").Value = 8296.6158
$ws.Range("M122").Value = -5846.6158
$ws.Range("H134").Value = 22863.445
$ws.Range("I134").Value = 31320.709
$ws.Range("J134").Value = 4136.643
$ws.Range("K134").Value = 93962.127
$ws.Range("L134").Value = 12409.929
$ws.Range("M134").Value = -91427.127
$ws.Range("N134").Value = -17479.929
$ws.Range("H136").Value = 1697.7333
$ws.Range("I136").Value = 1462.2142
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 4386.642599999999
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -1836.642599999999
$ws.Range("N136").Value = -20085

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H33").Value = 8443.083
$ws.Range("I33").Value = 85
$ws.Range("J33").Value = 33517.332
$ws.Range("K33").Value = 510
$ws.Range("L33").Value = 201103.992
$ws.Range("M33").Value = -227
$ws.Range("N33").Value = -201669.992
$ws.Range("H56").Value = 12506313
$ws.Range("I56").Value = 12506313
$ws.Range("K56").Value = 12506313
$ws.Range("M56").Value = -12505783
$ws.Range("H112").Value = 10054
$ws.Range("I112").Value = 10054
$ws.Range("K112").Value = 30162
$ws.Range("M112").Value = -29054

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H11").Value = 8180
$ws.Range("I11").Value = 900
$ws.Range("K11").Value = 900
$ws.Range("M11").Value = -761
$ws.Range("H57").Value = 8519.406
$ws.Range("I57").Value = 5777.2144
$ws.Range("J57").Value = 27714.75
$ws.Range("K57").Value = 5777.2144
$ws.Range("L57").Value = 27714.75
$ws.Range("M57").Value = -4957.2144
$ws.Range("N57").Value = -29354.75
$ws.Range("H97").Value = 834.5172
$ws.Range("J97").Value = 729
$ws.Range("L97").Value = 729
$ws.Range("N97").Value = -1721
$ws.Range("H122").Value = 59903.062
$ws.Range("I122").Value = 82061.26
$ws.Range("J122").Value = 3276.5557
$ws.Range("K122").Value = 246183.78
$ws.Range("L122").Value = 9829.667099999999
$ws.Range("M122").Value = -243733.78
$ws.Range("N122").Value = -14729.6671
$ws.Range("H126").Value = 3028.1333
$ws.Range("I126").Value = 2957.926
$ws.Range("J126").Value = 3660
$ws.Range("K126").Value = 8873.778
$ws.Range("L126").Value = 10980
$ws.Range("M126").Value = -6403.778
$ws.Range("N126").Value = -15920
$ws.Range("H132").Value = 4079.2222
$ws.Range("I132").Value = 3615
$ws.Range("J132").Value = 6748.5
$ws.Range("K132").Value = 10845
$ws.Range("L132").Value = 20245.5
$ws.Range("M132").Value = -8315
$ws.Range("N132").Value = -25305.5
$ws.Range("H137").Value = 70745
$ws.Range("J137").Value = 70745
$ws.Range("L137").Value = 70745
$ws.Range("N137").Value = -80945

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 11957.667
$ws.Range("I40").Value = 6000
$ws.Range("J40").Value = 14936.5
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 14936.5
$ws.Range("M40").Value = -5864
$ws.Range("N40").Value = -15208.5
$ws.Range("H41").Value = 42516.5
$ws.Range("I41").Value = 40033
$ws.Range("J41").Value = 45000
$ws.Range("K41").Value = 40033
$ws.Range("L41").Value = 45000
$ws.Range("M41").Value = -39595
$ws.Range("N41").Value = -45876
$ws.Range("H46").Value = 6824.92
$ws.Range("I46").Value = 5787.143
$ws.Range("J46").Value = 8145.727
$ws.Range("K46").Value = 5787.143
$ws.Range("L46").Value = 8145.727
$ws.Range("M46").Value = -5599.143
$ws.Range("N46").Value = -8521.726999999999
$ws.Range("H68").Value = 3166.6667
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 3166.6667
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -22488
$ws.Range("H93").Value = 2239.4412
$ws.Range("I93").Value = 2061.9
$ws.Range("K93").Value = 2061.9
$ws.Range("M93").Value = -813.9000000000001
$ws.Range("H100").Value = 66658.69
$ws.Range("I100").Value = 4454.3
$ws.Range("K100").Value = 4454.3
$ws.Range("M100").Value = -3913.3
$ws.Range("H122").Value = 5154.7
$ws.Range("I122").Value = 3878.75
$ws.Range("J122").Value = 7706.6
$ws.Range("K122").Value = 11636.25
$ws.Range("L122").Value = 23119.8
$ws.Range("M122").Value = -9186.25
$ws.Range("N122").Value = -28019.8

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H96").Value = 1978
$ws.Range("I96").Value = 1311
$ws.Range("K96").Value = 1311
$ws.Range("M96").Value = 62
$ws.Range("H106").Value = 54149.5
$ws.Range("J106").Value = 54149.5
$ws.Range("L106").Value = 54149.5
$ws.Range("N106").Value = -56673.5
$ws.Range("H126").Value = 2430.0952
$ws.Range("I126").Value = 2408.6875
$ws.Range("J126").Value = 2498.6
$ws.Range("K126").Value = 7226.0625
$ws.Range("L126").Value = 7495.799999999999
$ws.Range("M126").Value = -4756.0625
$ws.Range("N126").Value = -12435.8
$ws.Range("H132").Value = 208914.55
$ws.Range("I132").Value = 5148.2705
$ws.Range("J132").Value = 837193.94
$ws.Range("K132").Value = 15444.8115
$ws.Range("L132").Value = 2511581.82
$ws.Range("M132").Value = -12914.8115
$ws.Range("N132").Value = -2516641.82
$ws.Range("H136").Value = 1913.3636
$ws.Range("I136").Value = 1454.75
$ws.Range("J136").Value = 6499.5
$ws.Range("K136").Value = 4364.25
$ws.Range("L136").Value = 19498.5
$ws.Range("M136").Value = -1814.25
$ws.Range("N136").Value = -24598.5

Write-Host "Updated sheets: set=362 cleared=2"
